$d = $word.ActiveDocument

# 1/6/8. All three "May 22, 2022" dates move to "May 23, 2022"
$null = $d.Content.Find.Execute("May 22, 2022", $false, $false, $false, $false, $false, $true, 1, $false, "May 23, 2022", 2)

# 7. "July 21, 2022" -> "July 22, 2022"
$null = $d.Content.Find.Execute("July 21, 2022", $false, $false, $false, $false, $false, $true, 1, $false, "July 22, 2022", 2)

# 2. "...following allocution, entered the following sentence:" -> "...following allocution, the Court entered the following sentence:"
$null = $d.Content.Find.Execute("following allocution, entered the following sentence:", $false, $false, $false, $false, $false, $true, 1, $false, "following allocution, the Court entered the following sentence:", 2)

# 3. "Court costs are assessed for the highest degree charge in this case" -> "The Court ordered costs for the highest degree charge"
$null = $d.Content.Find.Execute("Court costs are assessed for the highest degree charge in this case", $false, $false, $false, $false, $false, $true, 1, $false, "The Court ordered costs for the highest degree charge", 2)

# 4. "Having been informed of the fines " -> "The Court informed Defendant of the fines "
$null = $d.Content.Find.Execute("Having been informed of the fines ", $false, $false, $false, $false, $false, $true, 1, $false, "The Court informed Defendant of the fines ", 2)

# 5. "owed, Defendant expressed an ability to pay " -> "owed, and Defendant expressed an ability to pay "
$null = $d.Content.Find.Execute("owed, Defendant expressed an ability to pay ", $false, $false, $false, $false, $false, $true, 1, $false, "owed, and Defendant expressed an ability to pay ", 2)
